$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Week1"

$ws.Range("A2").Value = 'DJ''s Quality Team'
$ws.Range("B2").Value = 120.0
$ws.Range("C2").Value = 59.0
$ws.Range("D2").Value = 61.0

$ws.Range("A3").Value = 'Samsquanches'
$ws.Range("B3").Value = 116.5
$ws.Range("C3").Value = 70.5
$ws.Range("D3").Value = 46.0

$ws.Range("A4").Value = 'Sánteros'
$ws.Range("B4").Value = 112.5
$ws.Range("C4").Value = 70.0
$ws.Range("D4").Value = 42.5

$ws.Range("A5").Value = 'Lundo’s Legends'
$ws.Range("B5").Value = 105.5
$ws.Range("C5").Value = 54.5
$ws.Range("D5").Value = 51.0

$ws.Range("A6").Value = 'SmokeWalkers'
$ws.Range("B6").Value = 101.5
$ws.Range("C6").Value = 64.0
$ws.Range("D6").Value = 37.5

$ws.Range("A7").Value = 'Swampnuts'
$ws.Range("B7").Value = 98.0
$ws.Range("C7").Value = 38.5
$ws.Range("D7").Value = 59.5

$ws.Range("A8").Value = 'GOD WILLS IT'
$ws.Range("B8").Value = 88.0
$ws.Range("C8").Value = 44.0
$ws.Range("D8").Value = 44.0

$ws.Range("A9").Value = 'MillerTime'
$ws.Range("B9").Value = 85.0
$ws.Range("C9").Value = 37.0
$ws.Range("D9").Value = 48.0

$ws.Range("A10").Value = 'PrimeTime'
$ws.Range("B10").Value = 83.0
$ws.Range("C10").Value = 51.5
$ws.Range("D10").Value = 31.5

$ws.Range("A11").Value = 'confusion'
$ws.Range("B11").Value = 81.0
$ws.Range("C11").Value = 47.0
$ws.Range("D11").Value = 34.0

$ws.Range("A12").Value = 'rainmaker'
$ws.Range("B12").Value = 72.0
$ws.Range("C12").Value = 28.5
$ws.Range("D12").Value = 43.5

$ws.Range("A13").Value = 'Epic7'
$ws.Range("B13").Value = 70.5
$ws.Range("C13").Value = 32.0
$ws.Range("D13").Value = 38.5

$ws.Range("A14").Value = 'KING JOSEPH 1 3 1'
$ws.Range("B14").Value = 67.5
$ws.Range("C14").Value = 16.0
$ws.Range("D14").Value = 51.5

$ws.Range("A15").Value = 'Splitfinger Skadoosh'
$ws.Range("B15").Value = 59.0
$ws.Range("C15").Value = 17.5
$ws.Range("D15").Value = 41.5
